$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 1.14
$ws.Range("N2").Value = 5.5
$ws.Range("S2").Value = 1.73
$ws.Range("T2").Value = 2.08
